# Didier's updates (part 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the old "Household (Optional)" rich-text label with "Disaggregates"
# wherever it appears in the level_lab column (H2:H4 all shared this string).
$ws.Range("H2:H4").Value = "Disaggregates"

# Update the active selection on the sheet to G5 (single cell), matching the
# saved sheetView state.
$ws.Range("G5").Select()
